$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.781.86"
$ws.Range('E2').Value = "'  -0.04%  "
$ws.Range('D3').Value = "'2.317.76"
$ws.Range('E3').Value = "'  +3.53%  "
$ws.Range('E4').Value = "'  -0.04%  "
$ws.Range('D5').Value = "'98.09"
$ws.Range('E5').Value = "'  +6.10%  "
$ws.Range('D6').Value = "'272.60"
$ws.Range('E6').Value = "'  +0.84%  "
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('E8').Value = "'  -0.05%  "
$ws.Range('D9').Value = "'0.627"
$ws.Range('E9').Value = "'  +0.60%  "
$ws.Range('D10').Value = "'45.48"
$ws.Range('E10').Value = "'  -1.61%  "
$ws.Range('D11').Value = "'0.0951"
$ws.Range('E11').Value = "'  -0.52%  "
$ws.Range('E12').Value = "'  -2.47%  "
$ws.Range('E13').Value = "'  -0.12%  "
$ws.Range('D14').Value = "'2.657.52"
$ws.Range('E14').Value = "'  +3.29%  "
$ws.Range('D15').Value = "'15.52"
$ws.Range('E15').Value = "'  +3.00%  "
$ws.Range('D16').Value = "'0.876"
$ws.Range('E16').Value = "'  +8.75%  "
$ws.Range('D17').Value = "'2.322.11"
$ws.Range('E17').Value = "'  +3.76%  "
$ws.Range('D18').Value = "'43.737.46"
$ws.Range('E18').Value = "'  -0.15%  "
$ws.Range('E19').Value = "'  +4.23%  "
$ws.Range('E20').Value = "'  +4.92%  "
$ws.Range('D21').Value = "'73.33"
$ws.Range('E21').Value = "'  +3.55%  "
$ws.Range('D22').Value = "'240.53"
$ws.Range('E22').Value = "'  +2.75%  "
$ws.Range('D23').Value = "'2.27"
$ws.Range('E23').Value = "'  -3.04%  "
$ws.Range('E24').Value = "'  +4.52%  "
$ws.Range('E25').Value = "'  -0.06%  "
$ws.Range('E26').Value = "'  +1.71%  "
$ws.Range('D27').Value = "'11.40"
$ws.Range('E27').Value = "'  -0.58%  "
$ws.Range('E28').Value = "'  -1.14%  "
$ws.Range('E29').Value = "'  +0.89%  "
$ws.Range('D30').Value = "'38.43"
$ws.Range('E30').Value = "'  -7.44%  "
$ws.Range('D31').Value = "'22.42"
$ws.Range('E31').Value = "'  +6.92%  "
$ws.Range('D32').Value = "'174.82"
$ws.Range('E33').Value = "'  -0.50%  "
$ws.Range('E34').Value = "'  -0.51%  "
$ws.Range('E35').Value = "'  +2.97%  "
$ws.Range('E36').Value = "'  -4.31%  "
$ws.Range('E37').Value = "'  +2.89%  "
$ws.Range('E38').Value = "'  +2.63%  "
$ws.Range('E39').Value = "'  -3.65%  "
$ws.Range('D40').Value = "'0.244"
$ws.Range('E40').Value = "'  +5.57%  "
$ws.Range('D41').Value = "'2.39"
$ws.Range('E41').Value = "'  +8.63%  "
$ws.Range('D42').Value = "'1.41"
$ws.Range('E42').Value = "'  +21.68%  "
$ws.Range('D43').Value = "'12.24"
$ws.Range('E43').Value = "'  -4.09%  "
$ws.Range('B44').Value = "'MultiversX"
$ws.Range('C44').Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range('D44').Value = "'63.11"
$ws.Range('E44').Value = "'  -0.87%  "
$ws.Range('B45').Value = "'FraxShare"
$ws.Range('C45').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('D45').Value = "'9.24"
$ws.Range('E45').Value = "'  +10.48%  "
$ws.Range('D46').Value = "'5.36"
$ws.Range('E46').Value = "'  -0.19%  "
$ws.Range('E47').Value = "'  +3.64%  "
$ws.Range('D48').Value = "'100.50"
$ws.Range('E48').Value = "'  +0.14%  "
$ws.Range('E49').Value = "'  +0.67%  "
$ws.Range('D50').Value = "'0.192"
$ws.Range('E50').Value = "'  +16.42%  "
$ws.Range('D51').Value = "'2.544.87"
$ws.Range('E51').Value = "'  +3.45%  "
